$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the containing path prefix from the fastq filenames in column F (rows 2-18)
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '^sequence/run_0674_samples/', ''
        $cell.Value = $newVal
    }
}

# Update selection to F2:F18
$ws.Range("F2:F18").Select()
